$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1733.3334
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 15
$ws.Range("H12").Value = 266.66666
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 250
$ws.Range("M12").Value = -130
$ws.Range("N12").Value = -590
$ws.Range("H29").Value = 2386.4666
$ws.Range("J29").Value = 2499.75
$ws.Range("L29").Value = 7499.25
$ws.Range("N29").Value = -8061.25
$ws.Range("H97").Value = 6946487
$ws.Range("J97").Value = 8335604
$ws.Range("L97").Value = 25006812
$ws.Range("N97").Value = -25007804
$ws.Range("H98").Value = 2110.923
$ws.Range("I98").Value = 2103.5833
$ws.Range("K98").Value = 2103.5833
$ws.Range("M98").Value = -605.5832999999998
$ws.Range("H122").Value = 2110.923
$ws.Range("I122").Value = 2103.5833
$ws.Range("K122").Value = 6310.749899999999
$ws.Range("M122").Value = -3860.749899999999
$ws.Range("H137").Value = 23812308
$ws.Range("I137").Value = 2600
$ws.Range("K137").Value = 7800
$ws.Range("M137").Value = -5250
$ws.Range("H138").Value = 5476.377
$ws.Range("I138").Value = 2553.1667
$ws.Range("K138").Value = 7659.500100000001
$ws.Range("M138").Value = -2519.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16708.5
$ws.Range("I32").Value = 16691.646
$ws.Range("K32").Value = 16691.646
$ws.Range("M32").Value = -16404.646
$ws.Range("H74").Value = 16669400
$ws.Range("I74").Value = 41668084
$ws.Range("J74").Value = 3610
$ws.Range("K74").Value = 41668084
$ws.Range("L74").Value = 3610
$ws.Range("M74").Value = -41667210
$ws.Range("N74").Value = -5358
$ws.Range("H77").Value = 16669400
$ws.Range("I77").Value = 41668084
$ws.Range("J77").Value = 3610
$ws.Range("K77").Value = 208340420
$ws.Range("L77").Value = 18050
$ws.Range("M77").Value = -208336052
$ws.Range("N77").Value = -26786
$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344
$ws.Range("H132").Value = 14837.77
$ws.Range("I132").Value = 17236.195
$ws.Range("J132").Value = 9441.3125
$ws.Range("K132").Value = 51708.585
$ws.Range("L132").Value = 28323.9375
$ws.Range("M132").Value = -49178.585
$ws.Range("N132").Value = -33383.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2680.8
$ws.Range("I86").Value = 2401.5
$ws.Range("J86").Value = 3099.75
$ws.Range("K86").Value = 2401.5
$ws.Range("L86").Value = 3099.75
$ws.Range("M86").Value = -1278.5
$ws.Range("N86").Value = -5345.75
$ws.Range("H89").Value = 2680.8
$ws.Range("I89").Value = 2401.5
$ws.Range("J89").Value = 3099.75
$ws.Range("K89").Value = 12007.5
$ws.Range("L89").Value = 15498.75
$ws.Range("M89").Value = -6391.5
$ws.Range("N89").Value = -26730.75
$ws.Range("H99").Value = 63803864
$ws.Range("I99").Value = 85071040
$ws.Range("K99").Value = 85071040
$ws.Range("M99").Value = -85069542

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 268.53845
$ws.Range("I7").Value = 188.71428
$ws.Range("J7").Value = 361.66666
$ws.Range("K7").Value = 188.71428
$ws.Range("L7").Value = 361.66666
$ws.Range("M7").Value = -75.71428
$ws.Range("N7").Value = -587.66666
$ws.Range("H58").Value = 1249.44
$ws.Range("I58").Value = 1118.3478
$ws.Range("K58").Value = 1118.3478
$ws.Range("M58").Value = -915.3478
$ws.Range("H70").Value = 100000
$ws.Range("J70").Value = 100000
$ws.Range("L70").Value = 100000
$ws.Range("N70").Value = -100630
$ws.Range("H73").Value = 100000
$ws.Range("J73").Value = 100000
$ws.Range("L73").Value = 100000
$ws.Range("N73").Value = -102184
$ws.Range("H86").Value = 10721.342
$ws.Range("I86").Value = 10802.723
$ws.Range("J86").Value = 10648.1
$ws.Range("K86").Value = 10802.723
$ws.Range("L86").Value = 10648.1
$ws.Range("M86").Value = -9679.723
$ws.Range("N86").Value = -12894.1
$ws.Range("H89").Value = 10721.342
$ws.Range("I89").Value = 10802.723
$ws.Range("J89").Value = 10648.1
$ws.Range("K89").Value = 54013.615
$ws.Range("L89").Value = 53240.5
$ws.Range("M89").Value = -48397.615
$ws.Range("N89").Value = -64472.5
$ws.Range("H132").Value = 27790828
$ws.Range("I132").Value = 32529750
$ws.Range("K132").Value = 97589250
$ws.Range("M132").Value = -97586720
$ws.Range("H134").Value = 1675.7179
$ws.Range("I134").Value = 1633.6487
$ws.Range("K134").Value = 4900.9461
$ws.Range("M134").Value = -2365.9461
$ws.Range("H136").Value = 1249.44
$ws.Range("I136").Value = 1118.3478
$ws.Range("K136").Value = 3355.0434
$ws.Range("M136").Value = -805.0434
$ws.Range("H141").Value = 141022.9
$ws.Range("J141").Value = 141022.9
$ws.Range("L141").Value = 141022.9
$ws.Range("N141").Value = -151382.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 299.5
$ws.Range("I14").Value = 299.5
$ws.Range("K14").Value = 898.5
$ws.Range("M14").Value = -725.5
$ws.Range("H42").Value = 18500
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 18500
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 55500
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = -56568
$ws.Range("H140").Value = 7748.0435
$ws.Range("I140").Value = 2950.5
$ws.Range("J140").Value = 18713.857
$ws.Range("K140").Value = 8851.5
$ws.Range("L140").Value = 56141.571
$ws.Range("M140").Value = -3671.5
$ws.Range("N140").Value = -66501.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H80").Value = 987922.75
$ws.Range("J80").Value = 22100.834
$ws.Range("L80").Value = 22100.834
$ws.Range("N80").Value = -24096.834
$ws.Range("H83").Value = 987922.75
$ws.Range("J83").Value = 22100.834
$ws.Range("L83").Value = 110504.17
$ws.Range("N83").Value = -120488.17
$ws.Range("H113").Value = 3456.7273
$ws.Range("J113").Value = 3699.8333
$ws.Range("L113").Value = 3699.8333
$ws.Range("N113").Value = -8039.8333
$ws.Range("H132").Value = 6792.28
$ws.Range("I132").Value = 6488.4287
$ws.Range("K132").Value = 19465.2861
$ws.Range("M132").Value = -16935.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4573.074
$ws.Range("I46").Value = 964.6667
$ws.Range("K46").Value = 964.6667
$ws.Range("M46").Value = -776.6667
$ws.Range("H61").Value = 2450
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("H68").Value = 4551854
$ws.Range("I68").Value = 22727272
$ws.Range("J68").Value = 7999.5
$ws.Range("K68").Value = 22727272
$ws.Range("L68").Value = 7999.5
$ws.Range("M68").Value = -22726523
$ws.Range("N68").Value = -9497.5
$ws.Range("H71").Value = 4551854
$ws.Range("I71").Value = 22727272
$ws.Range("J71").Value = 7999.5
$ws.Range("K71").Value = 113636360
$ws.Range("L71").Value = 39997.5
$ws.Range("M71").Value = -113632616
$ws.Range("N71").Value = -47485.5
$ws.Range("H82").Value = 6251530
$ws.Range("I82").Value = 10417083
$ws.Range("J82").Value = 3200
$ws.Range("K82").Value = 10417083
$ws.Range("L82").Value = 3200
$ws.Range("M82").Value = -10416722
$ws.Range("N82").Value = -3922
$ws.Range("H85").Value = 6251530
$ws.Range("I85").Value = 10417083
$ws.Range("J85").Value = 3200
$ws.Range("K85").Value = 10417083
$ws.Range("L85").Value = 3200
$ws.Range("M85").Value = -10415835
$ws.Range("N85").Value = -5696
$ws.Range("H113").Value = 2450
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H122").Value = 15371
$ws.Range("I122").Value = 21997.5
$ws.Range("J122").Value = 13162.167
$ws.Range("K122").Value = 65992.5
$ws.Range("L122").Value = 39486.501
$ws.Range("M122").Value = -63542.5
$ws.Range("N122").Value = -44386.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17551492
$ws.Range("I62").Value = 26317238
$ws.Range("K62").Value = 26317238
$ws.Range("M62").Value = -26316614
$ws.Range("H65").Value = 17551492
$ws.Range("I65").Value = 26317238
$ws.Range("K65").Value = 131586190
$ws.Range("M65").Value = -131583070
$ws.Range("H107").Value = 419
$ws.Range("J107").Value = 374.8
$ws.Range("L107").Value = 1124.4
$ws.Range("N107").Value = -4964.4
$ws.Range("H113").Value = 1046.2222
$ws.Range("I113").Value = 861.3
$ws.Range("K113").Value = 2583.9
$ws.Range("M113").Value = -413.8999999999996
$ws.Range("H132").Value = 20001330
$ws.Range("I132").Value = 877.9286
$ws.Range("J132").Value = 45456452
$ws.Range("K132").Value = 2633.7858
$ws.Range("L132").Value = 136369356
$ws.Range("M132").Value = -103.7857999999997
$ws.Range("N132").Value = -136374416
$ws.Range("H136").Value = 8582.909
$ws.Range("I136").Value = 4141.9487
$ws.Range("K136").Value = 12425.8461
$ws.Range("M136").Value = -9875.846099999999
